$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197, shifting existing rows 197:216 down to 198:217.
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new data record.
$ws.Cells.Item(197, 1).Value = 3
$ws.Cells.Item(197, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(197, 3).Value = "Coquimbo"
$ws.Cells.Item(197, 4).Value = 44578
$ws.Cells.Item(197, 4).NumberFormat = $ws.Cells.Item(198, 4).NumberFormat
$ws.Cells.Item(197, 5).Value = 5
$ws.Cells.Item(197, 6).Value = 100112001
$ws.Cells.Item(197, 7).Value = "Berenjena"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 70
$ws.Cells.Item(197, 11).Value = 9000
$ws.Cells.Item(197, 12).Value = 9500
$ws.Cells.Item(197, 13).Value = 9214
$ws.Cells.Item(197, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(197, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(197, 16).Value = 154
$ws.Cells.Item(197, 17).Value = 60
$ws.Cells.Item(197, 18).Value = "Hortaliza"
